# Add two new market test-data sheets (Italy, Spain) based on the existing
# "Norway" sheet template, matching the layout used by the other market tabs.

$wb = $excel.ActiveWorkbook
$norway = $wb.Worksheets.Item("Norway")

# --- Spain sheet -----------------------------------------------------
# Copy Norway right after itself to keep all formatting/merged cells/column
# widths intact, then rename and fill in the market-specific values.
$norway.Copy($null, $norway)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3442/T2125"
$null = $spain.Range("A6").Select()

# --- Italy sheet -------------------------------------------------------
# Copy Norway again (also placed right after Norway), which pushes it in
# between Norway and the just-created Spain sheet.
$norway.Copy($null, $norway)
$italy = $wb.Worksheets.Item(6)
$italy.Name = "Italy"
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3443/T1972"

# Italy becomes the active/selected tab.
$null = $italy.Activate()
$null = $italy.Range("A8").Select()
